$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout (before edit):
#  Row2 (s=1): TWONTO | Super_Class
#  Row3 (s=2): instrumentation | Meter
#  Row4: air_duct_segment | Silencer
#  Row5: cable_segment | Electrical Power Line
#  Row6: electrical_panel_or_cabinet | BUS bar greater than 750 volts...
#  Row7: instrument_gauge_or_display | Pressure Indicator
#
# Target layout (after edit):
#  Row1 (s=3, wrapped, tall): <SPARQL query text>
#  Row2 (s=1): TWONTO | Super_Class
#  Row3 (s=2): instrument gauge or display | Pressure Indicator
#  Row4: cable segment | Electrical Power Line
#  Row5: air duct segment | Silencer
#
# Row 1 is currently empty/unused, so the new SPARQL row can be written
# straight into it - no row insertion/shifting required.

$sparql = "PREFIX rdf: <http://www.w3.org/1999/02/22-rdf-syntax-ns#>`nPREFIX owl: <http://www.w3.org/2002/07/owl#>`nPREFIX rdfs: <http://www.w3.org/2000/01/rdf-schema#>`nPREFIX xsd: <http://www.w3.org/2001/XMLSchema#>`nPREFIX tw: <http://www.toronto.ca/TWONTO#>`nSELECT (STR(?label) as ?TWONTO) (STR(?object) as ?Avantis)`nWHERE { `n    ?entityIRI tw:is_superclass_of_avantis_class ?object ;`n              rdfs:label ?label .`n}"

$ws.Cells.Item(1, 1).Value = $sparql
$ws.Cells.Item(1, 1).WrapText = $true
$ws.Rows.Item(1).RowHeight = 159.5

# Overwrite values in place (this does not disturb existing cell formatting):
#  Row3 already carries the banded style (s=2) - reuse it for the surviving
#  "instrument gauge or display"/"Pressure Indicator" row.
$ws.Cells.Item(3, 1).Value = "instrument gauge or display"
$ws.Cells.Item(3, 2).Value = "Pressure Indicator"
#  Row4/Row5 already carry no special style - reuse them for the other two
#  surviving (unstyled) rows.
$ws.Cells.Item(4, 1).Value = "cable segment"
$ws.Cells.Item(4, 2).Value = "Electrical Power Line"
$ws.Cells.Item(5, 1).Value = "air duct segment"
$ws.Cells.Item(5, 2).Value = "Silencer"

# Drop the two now-superfluous trailing rows (their content has already been
# relocated above). Delete the higher-numbered row first.
$ws.Rows.Item(7).Delete() | Out-Null
$ws.Rows.Item(6).Delete() | Out-Null

# Resize Table1 to its new extent (A2:B7 -> A2:B5).
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A2:B5"))

$ws.Range("B7").Select() | Out-Null
